$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Projektplan")
$ws.Rows("18:18").Delete()
